$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 692. This pushes the current rows 692..782 down to 693..783,
# preserving all of their data untouched (shift-down), exactly like Excel's own
# "Insert Row" behaviour.
$ws.Rows(692).Insert()

# Populate the newly-inserted row 692 with the new week's record (a new Cebollin
# price entry for Femacal de La Calera). The non-varying columns match every
# other row in this data block.
$ws.Cells.Item(692, 1).Value = 3
$ws.Cells.Item(692, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(692, 3).Value = "Coquimbo"
$ws.Cells.Item(692, 4).Value = 45124
$ws.Cells.Item(692, 5).Value = 5
$ws.Cells.Item(692, 6).Value = 100112037
$ws.Cells.Item(692, 7).Value = "Cebollín"
$ws.Cells.Item(692, 8).Value = "Sin especificar"
$ws.Cells.Item(692, 9).Value = "Primera"
$ws.Cells.Item(692, 10).Value = 240
$ws.Cells.Item(692, 11).Value = 4000
$ws.Cells.Item(692, 12).Value = 4300
$ws.Cells.Item(692, 13).Value = 4138
$ws.Cells.Item(692, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(692, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(692, 16).Value = 115
$ws.Cells.Item(692, 17).Value = 36
$ws.Cells.Item(692, 18).Value = "Hortaliza"
